$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- FilesTab query (column B, row 4) was rewritten to add file-size
#     formatting (Bytes/KB/MB/GB/TB), a Format column, and a Sample ID
#     lookup via an OPTIONAL MATCH back through the sample node.
$newFilesQuery = @"
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Beagle','Mixed Breed']and diag.disease_term in ['Bladder Cancer','Healthy Control'] and diag.primary_disease_site in ['Bladder, Urethra']
WITH DISTINCT f, parent, c, demo, diag, s
OPTIONAL MATCH (f)-[*]->(samp:sample)
OPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp)
WITH
        f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value, 
        10^precision AS factor,
        units[i] as unit
WITH    
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN coalesce(f.file_name, '') AS ``File Name``, 
 coalesce(f.file_format, '') AS ``Format``,
        coalesce(f.file_type, '') AS ``File Type``, 
      CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(labels(parent)[0], '') AS ``Association``,
        coalesce(f.file_description, '') AS ``Description``,
   coalesce(samp.sample_id, '') AS ``Sample ID``,
        coalesce(c.case_id, '') AS ``Case ID``, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis 
        Order By f.file_name LIMIT 100
"@

$ws.Range("B4").Value = $newFilesQuery

# --- Row heights re-settled (author resized the window / re-wrapped text)
$ws.Range("A2").RowHeight = 348
$ws.Range("A3").RowHeight = 275.5
$ws.Range("A4").RowHeight = 409.5

# --- Selection / scroll moved up one row after the edit
$ws.Range("C4").Select() | Out-Null
